$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.85
$ws.Range("H2").Value = 3.75
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 2.5
$ws.Range("K2").Value = 2.4
$ws.Range("L2").Value = 4
$ws.Range("O2").Value = 1.18
$ws.Range("P2").Value = 4.5
$ws.Range("S2").Value = 1.29
$ws.Range("U2").Value = 1.53
$ws.Range("V2").Value = 2.38
$ws.Range("W2").Value = 10
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 13
$ws.Range("AD2").Value = 7.5
$ws.Range("AG2").Value = 126
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 13
$ws.Range("AL2").Value = 26
$ws.Range("AM2").Value = 29
$ws.Range("AO2").Value = 9.5
$ws.Range("AQ2").Value = 29
$ws.Range("AU2").Value = 7.5
$ws.Range("AX2").Value = 6
$ws.Range("AY2").Value = 19
$ws.Range("AZ2").Value = 23
# Row 3
$ws.Range("G3").Value = 2.55
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 2.62
$ws.Range("N3").Value = 10.9
$ws.Range("O3").Value = 1.21
$ws.Range("P3").Value = 3.95
$ws.Range("Q3").Value = 1.62
$ws.Range("R3").Value = 2.02
$ws.Range("S3").Value = 1.33
$ws.Range("T3").Value = 3.04
$ws.Range("U3").Value = 1.54
$ws.Range("V3").Value = 2.32
$ws.Range("AB3").Value = 17
$ws.Range("AD3").Value = 5.4
$ws.Range("AE3").Value = 8.75
$ws.Range("AF3").Value = 27
$ws.Range("AH3").Value = 9
$ws.Range("AI3").Value = 13
$ws.Range("AK3").Value = 26
$ws.Range("AT3").Value = 2.85
$ws.Range("AU3").Value = 6
$ws.Range("AX3").Value = 4.8
# Row 4
$ws.Range("L4").Value = 4.25
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 9.449999999999999
$ws.Range("O4").Value = 1.2
$ws.Range("P4").Value = 4.05
$ws.Range("R4").Value = 2.02
$ws.Range("S4").Value = 1.31
$ws.Range("T4").Value = 3.15
$ws.Range("U4").Value = 1.68
$ws.Range("V4").Value = 2.14
$ws.Range("W4").Value = 7.1
$ws.Range("X4").Value = 7.5
$ws.Range("Y4").Value = 7
$ws.Range("AB4").Value = 18
$ws.Range("AE4").Value = 12
$ws.Range("AF4").Value = 45
$ws.Range("AM4").Value = 29
$ws.Range("AS4").Value = 200
$ws.Range("AV4").Value = 60
$ws.Range("AX4").Value = 6
$ws.Range("BB4").Value = 150
# Row 5
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 15
# Row 6
$ws.Range("G6").Value = 1.65
$ws.Range("I6").Value = 4.85
$ws.Range("J6").Value = 2.2
$ws.Range("K6").Value = 2.15
$ws.Range("L6").Value = 5
$ws.Range("O6").Value = 1.26
$ws.Range("P6").Value = 3.15
$ws.Range("Q6").Value = 1.78
$ws.Range("R6").Value = 1.83
$ws.Range("S6").Value = 1.35
$ws.Range("T6").Value = 2.94
$ws.Range("U6").Value = 1.78
$ws.Range("V6").Value = 1.83
$ws.Range("W6").Value = 6.9
$ws.Range("X6").Value = 7.7
$ws.Range("Z6").Value = 12.5
$ws.Range("AB6").Value = 26
$ws.Range("AC6").Value = 10.5
$ws.Range("AD6").Value = 7.1
$ws.Range("AE6").Value = 16
$ws.Range("AF6").Value = 75
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 13.5
$ws.Range("AI6").Value = 29
$ws.Range("AJ6").Value = 15.5
$ws.Range("AK6").Value = 90
$ws.Range("AL6").Value = 50
$ws.Range("AM6").Value = 50
$ws.Range("AN6").Value = 3.45
$ws.Range("AO6").Value = 7.9
$ws.Range("AP6").Value = 17
$ws.Range("AQ6").Value = 26
$ws.Range("AS6").Value = 250
$ws.Range("AT6").Value = 2.72
$ws.Range("AU6").Value = 7.4
$ws.Range("AV6").Value = 70
$ws.Range("AX6").Value = 6.5
$ws.Range("AY6").Value = 28
$ws.Range("AZ6").Value = 32
$ws.Range("BA6").Value = 175
$ws.Range("BB6").Value = 200
$ws.Range("BC6").Value = 450
# Row 7
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 5.2
$ws.Range("P7").Value = 2.8
$ws.Range("W7").Value = 5.8
$ws.Range("AA7").Value = 14.5
$ws.Range("AE7").Value = 17.5
$ws.Range("AJ7").Value = 16.5
$ws.Range("AL7").Value = 55
$ws.Range("AM7").Value = 60
$ws.Range("AP7").Value = 17.5
$ws.Range("AT7").Value = 2.55
$ws.Range("AU7").Value = 7.5
$ws.Range("AV7").Value = 70
$ws.Range("AY7").Value = 29
$ws.Range("AZ7").Value = 32
$ws.Range("BA7").Value = 175
$ws.Range("BB7").Value = 200
$ws.Range("BC7").Value = 450
# Row 8
$ws.Range("I8").Value = 5.5
$ws.Range("J8").Value = 2.05
$ws.Range("L8").Value = 5.2
$ws.Range("Q8").Value = 1.57
$ws.Range("R8").Value = 2.12
$ws.Range("U8").Value = 1.62
$ws.Range("Y8").Value = 7.8
$ws.Range("AC8").Value = 13.5
$ws.Range("AD8").Value = 7.8
$ws.Range("AE8").Value = 14
$ws.Range("AH8").Value = 18.5
$ws.Range("AJ8").Value = 17
$ws.Range("AO8").Value = 7.3
$ws.Range("AX8").Value = 7.2
$ws.Range("AY8").Value = 29
$ws.Range("AZ8").Value = 28
$ws.Range("BA8").Value = 175
$ws.Range("BB8").Value = 175
$ws.Range("BC8").Value = 300
# Row 9
$ws.Range("G9").Value = 1.85
$ws.Range("H9").Value = 3.5
$ws.Range("I9").Value = 4.5
$ws.Range("J9").Value = 2.5
$ws.Range("K9").Value = 2.1
$ws.Range("L9").Value = 5
$ws.Range("X9").Value = 8
$ws.Range("Y9").Value = 8.5
$ws.Range("Z9").Value = 15
$ws.Range("AD9").Value = 6.5
$ws.Range("AE9").Value = 17
$ws.Range("AH9").Value = 11
$ws.Range("AI9").Value = 21
$ws.Range("AJ9").Value = 15
$ws.Range("AK9").Value = 51
$ws.Range("AL9").Value = 41
$ws.Range("AN9").Value = 3.75
$ws.Range("AO9").Value = 10
$ws.Range("AQ9").Value = 34
$ws.Range("AR9").Value = 51
$ws.Range("AX9").Value = 6
$ws.Range("AY9").Value = 26
$ws.Range("BA9").Value = 81
$ws.Range("BB9").Value = 126
$ws.Range("BC9").Value = 301
# Row 10
$ws.Range("Q10").Value = 1.7
$ws.Range("R10").Value = 2.1
# Row 11
$ws.Range("G11").Value = 4.8
$ws.Range("I11").Value = 1.65
$ws.Range("J11").Value = 5.2
$ws.Range("K11").Value = 2.12
$ws.Range("L11").Value = 2.22
$ws.Range("P11").Value = 2.9
$ws.Range("S11").Value = 1.44
$ws.Range("T11").Value = 2.62
$ws.Range("U11").Value = 2.05
$ws.Range("W11").Value = 11.5
$ws.Range("X11").Value = 27
$ws.Range("Y11").Value = 16.5
$ws.Range("Z11").Value = 90
$ws.Range("AA11").Value = 55
$ws.Range("AB11").Value = 65
$ws.Range("AH11").Value = 5.8
$ws.Range("AI11").Value = 7
$ws.Range("AJ11").Value = 8.5
$ws.Range("AK11").Value = 12
$ws.Range("AL11").Value = 14.5
$ws.Range("AN11").Value = 6.5
$ws.Range("AO11").Value = 29
$ws.Range("AP11").Value = 37
$ws.Range("AT11").Value = 2.62
$ws.Range("AU11").Value = 8.25
$ws.Range("AV11").Value = 90
$ws.Range("AY11").Value = 8.25
$ws.Range("AZ11").Value = 19.5
$ws.Range("BA11").Value = 28
$ws.Range("BB11").Value = 70
